# Fulfillment for General Questions
# Updates the "intents-en" sheet: adds training phrases for age/name/sex
# intents, switches the parameter placeholder syntax to the "$name"-style
# references, relocates the fallback_age_set / fallback_name_set /
# fallback_sex_set rows (now with response phrases) up next to their
# corresponding *_set intents, and moves initial_symptom_set down below them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("intents-en")

# --- Row 4: privacy_policy_yes (unchanged content, shared-string index shift only) ---
$ws.Range("A4").Value = 'privacy_policy_yes'
$ws.Range("B4").Value = "en"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 'PRIVACY_POLICY'
$ws.Range("G4").Value = 'PRIVACY-POLICY-YES'
$ws.Range("H4").Value = 'Yes|Okay I will|Why not|Yes that''s alright|Yes I do|Exactly|Of course|Yep that''s ok|Okay|Ok'

# --- Row 5: privacy_policy_no ---
$ws.Range("A5").Value = 'privacy_policy_no'
$ws.Range("B5").Value = "en"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 'PRIVACY_POLICY'
$ws.Range("G5").Value = 'PRIVACY-POLICY-NO'
$ws.Range("H5").Value = 'No way|No|Na|I can''t|No I cannot|Don''t|Nope|I disagree|Of course not|No thanks|Not right now|Nah'

# --- Row 6: language_set -- parameter placeholder now uses "$language" ---
$ws.Range("A6").Value = 'language_set'
$ws.Range("B6").Value = "en"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 'LANGUAGE'
$ws.Range("G6").Value = 'LANGUAGE-SET'
$ws.Range("H6").Value = '{@language:english}|I want {@language:tagalog}|I choose {@language:filipino}|{@language:english} please|Let''s go with {@language:tagalog}|Change it to {@language:filipino}|Use {@language:filipino} instead'
$ws.Range("J6").Value = '{@language:$language:1:0}'

# --- Row 7: language_change -- parameter placeholder now uses "$language" ---
$ws.Range("A7").Value = 'language_change'
$ws.Range("B7").Value = "en"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 0
$ws.Range("G7").Value = 'LANGUAGE-CHANGE'
$ws.Range("H7").Value = 'I want to change my language|Let me change language|Change language|Change my language to {@language:english}|Use {@language:english} instead'
$ws.Range("J7").Value = '{@language:$language:0:0}'
$ws.Range("L7").Value = ' '

# --- Row 8: fallback_privacy_policy ---
$ws.Range("A8").Value = 'fallback_privacy_policy'
$ws.Range("B8").Value = "en"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 'PRIVACY_POLICY'
$ws.Range("I8").Value = 'I''m sorry, please answer the question. Do you agree to my terms of service?|Sorry? Do you agree to my terms of service?'

# --- Row 9: fallback_language_set ---
$ws.Range("A9").Value = 'fallback_language_set'
$ws.Range("B9").Value = "en"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 'LANGUAGE'
$ws.Range("I9").Value = 'What language?|Sorry, I did not quite get that. What language do you prefer?'

# --- Row 10: general -- no longer carries an inputContext ---
$ws.Range("A10").Value = 'general'
$ws.Range("B10").Value = "en"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 0
$ws.Range("E10").Clear()
$ws.Range("G10").Value = 'GENERAL'

# --- Row 11: age_set -- add training phrases, "$age" parameter syntax ---
$ws.Range("A11").Value = 'age_set'
$ws.Range("B11").Value = "en"
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 'AGE'
$ws.Range("G11").Value = 'AGE'
$ws.Range("H11").Value = 'I am {@sys.age:20 years old}|{@sys.age:30 years old}|I''m {@sys.age:30 years old} right now|{@sys.age:25 old}|{@sys.age:18 years}'
$ws.Range("J11").Value = '{@sys.age:$age:1:0}'

# --- Row 12: name_set -- add training phrases, "$name" parameter syntax ---
$ws.Range("A12").Value = 'name_set'
$ws.Range("B12").Value = "en"
$ws.Range("C12").Value = 2
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 'NAME'
$ws.Range("G12").Value = 'NAME'
$ws.Range("H12").Value = '{@sys.any:Steven}|You can call me {@sys.any:Mark}|Address me as {@sys.any:Robert}|My name is {@sys.any:Maxson}|Call me {@sys.any:Roy}'
$ws.Range("J12").Value = '{@sys.any:$name:1:0}'

# --- Row 13: sex_set -- add training phrases, "$sex" parameter syntax ---
$ws.Range("A13").Value = 'sex_set'
$ws.Range("B13").Value = "en"
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 'SEX'
$ws.Range("G13").Value = 'SEX'
$ws.Range("H13").Value = '{@sex:male}|I am a {@sex:girl}|I''m a {@sex:man}|My sex is {@sex:female}|My gender is {@sex:male}'
$ws.Range("J13").Value = '{@sex:$sex:1:0}'

# --- Row 14: fallback_age_set (relocated here, now with a response phrase) ---
$ws.Range("A14").Value = 'fallback_age_set'
$ws.Range("B14").Value = "en"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 'AGE'
$ws.Range("G14").Clear()
$ws.Range("I14").Value = 'Again how old are you?|I''m sorry, how old?|I did not get that, what is your age?'

# --- Row 15: fallback_name_set (relocated here, now with a response phrase) ---
$ws.Range("A15").Value = 'fallback_name_set'
$ws.Range("B15").Value = "en"
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 'NAME'
$ws.Range("I15").Value = 'What name?'

# --- Row 16: fallback_sex_set (relocated here, now with a response phrase) ---
$ws.Range("A16").Value = 'fallback_sex_set'
$ws.Range("B16").Value = "en"
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 'SEX'
$ws.Range("I16").Value = 'What sex?|What is your biological sex?'

# --- Row 17: initial_symptom_set (moved down from row 14) ---
$ws.Range("A17").Value = 'initial_symptom_set'
$ws.Range("B17").Value = "en"
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 'INITIAL'
$ws.Range("G17").Value = 'INITIAL'

# --- Rows 19-21: vacated by the relocated fallback_* rows above; only the
#     B:D formatting remnants stay behind. ---
$ws.Range("A19").Clear()
$ws.Range("E19").Clear()
$ws.Range("B19:D19").ClearContents()
$ws.Range("B19:D19").HorizontalAlignment = 1

$ws.Range("A20").Clear()
$ws.Range("E20").Clear()
$ws.Range("B20:D20").ClearContents()
$ws.Range("B20:D20").HorizontalAlignment = 1

$ws.Range("A21").Clear()
$ws.Range("E21").Clear()
$ws.Range("B21:D21").ClearContents()
$ws.Range("B21:D21").HorizontalAlignment = 1

# --- Row heights for the rows whose wrapped text now spans more lines ---
$ws.Rows.Item(11).RowHeight = 56.25
$ws.Rows.Item(12).RowHeight = 56.25
$ws.Rows.Item(13).RowHeight = 45
$ws.Rows.Item(14).RowHeight = 33.75
$ws.Rows.Item(16).RowHeight = 22.5

# --- Restore the view/selection to where the edit left off ---
$ws.Activate()
$ws.Range("G12").Select()
